$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# A leading apostrophe forces Excel to store the value as literal text
# (matching the source workbook's inlineStr cells) instead of coercing
# numeric-looking strings (e.g. "7.50", "0.518") into numbers.


# Row 2
$ws.Range("D2").Value = "'67.813.71"
$ws.Range("E2").Value = "'  +1.39%  "

# Row 3
$ws.Range("D3").Value = "'2.504.44"
$ws.Range("E3").Value = "'  +1.86%  "

# Row 4
$ws.Range("E4").Value = "'  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'588.12"
$ws.Range("E5").Value = "'  +1.19%  "

# Row 6
$ws.Range("D6").Value = "'175.94"
$ws.Range("E6").Value = "'  +4.20%  "

# Row 7
$ws.Range("E7").Value = "'  -0.08%  "

# Row 8
$ws.Range("E8").Value = "'  +1.45%  "

# Row 9
$ws.Range("D9").Value = "'0.142"
$ws.Range("E9").Value = "'  +6.16%  "

# Row 11
$ws.Range("E11").Value = "'  +4.34%  "

# Row 12
$ws.Range("E12").Value = "'  +1.50%  "

# Row 13
$ws.Range("D13").Value = "'2.943.82"

# Row 14
$ws.Range("D14").Value = "'25.84"
$ws.Range("E14").Value = "'  +2.90%  "

# Row 15
$ws.Range("D15").Value = "'67.596.53"
$ws.Range("E15").Value = "'  +1.18%  "

# Row 16
$ws.Range("D16").Value = "'0.0000173"
$ws.Range("E16").Value = "'  +3.19%  "

# Row 17
$ws.Range("D17").Value = "'2.496.58"
$ws.Range("E17").Value = "'  +1.56%  "

# Row 18
$ws.Range("D18").Value = "'11.13"
$ws.Range("E18").Value = "'  +2.38%  "

# Row 19
$ws.Range("D19").Value = "'7.50"
$ws.Range("E19").Value = "'  +2.48%  "

# Row 20
$ws.Range("D20").Value = "'352.55"
$ws.Range("E20").Value = "'  +1.73%  "

# Row 21
$ws.Range("E21").Value = "'  +2.58%  "

# Row 22
$ws.Range("E22").Value = "'  +0.05%  "

# Row 23
$ws.Range("D23").Value = "'70.71"
$ws.Range("E23").Value = "'  +3.34%  "

# Row 24
$ws.Range("D24").Value = "'4.27"
$ws.Range("E24").Value = "'  +2.31%  "

# Row 25
$ws.Range("E25").Value = "'  -0.15%  "

# Row 26
$ws.Range("E26").Value = "'  +2.17%  "

# Row 27
$ws.Range("D27").Value = "'2.628.06"
$ws.Range("E27").Value = "'  +1.73%  "

# Row 28
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "'  -0.30%  "

# Row 29
$ws.Range("D29").Value = "'0.0₃0915"
$ws.Range("E29").Value = "'  +2.91%  "

# Row 30
$ws.Range("D30").Value = "'513.61"
$ws.Range("E30").Value = "'  +0.99%  "

# Row 31
$ws.Range("D31").Value = "'7.87"
$ws.Range("E31").Value = "'  +4.23%  "

# Row 32
$ws.Range("E32").Value = "'  +4.59%  "

# Row 33
$ws.Range("E33").Value = "'  +2.13%  "

# Row 34
$ws.Range("E34").Value = "'  +0.10%  "

# Row 35
$ws.Range("E35").Value = "'  +8.46%  "

# Row 36
$ws.Range("D36").Value = "'161.96"
$ws.Range("E36").Value = "'  +2.25%  "

# Row 37
$ws.Range("D37").Value = "'18.50"
$ws.Range("E37").Value = "'  +2.24%  "

# Row 38
$ws.Range("E38").Value = "'  +0.34%  "

# Row 40
$ws.Range("D40").Value = "'1.76"
$ws.Range("E40").Value = "'  +6.51%  "

# Row 42
$ws.Range("D42").Value = "'0.332"
$ws.Range("E42").Value = "'  +3.02%  "

# Row 43
$ws.Range("E43").Value = "'  +3.18%  "

# Row 44
$ws.Range("E44").Value = "'  +4.48%  "

# Row 45
$ws.Range("D45").Value = "'145.57"
$ws.Range("E45").Value = "'  +3.52%  "

# Row 46
$ws.Range("E46").Value = "'  +3.42%  "

# Row 47
$ws.Range("B47").Value = "'ARBITRUM"
$ws.Range("C47").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'0.518"
$ws.Range("E47").Value = "'  +2.55%  "

# Row 48
$ws.Range("B48").Value = "'BabyDogeCoin"
$ws.Range("C48").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.0₆0258"
$ws.Range("E48").Value = "'  +2.24%  "

# Row 49
$ws.Range("E49").Value = "'  +2.89%  "

# Row 50
$ws.Range("E50").Value = "'  +3.04%  "

# Row 51
$ws.Range("D51").Value = "'0.588"
$ws.Range("E51").Value = "'  +1.65%  "
